$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill the Season column (B) for each seasonal block with the season label,
# lower-casing "winter"/"summer" while leaving "Monsoon"/"Post-monsoon" as-is.
$ws.Range("B3:B10").Value = "winter"
$ws.Range("B11:B16").Value = "summer"
$ws.Range("B17:B25").Value = "Monsoon"
$ws.Range("B26:B30").Value = "Post-monsoon"

# Update the selected cell to match the saved view state.
$ws.Range("B32").Select()
